# completed tank titrations 0413
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Append the new titration data row (row 81), continuing the shared
# formula series that computes "% off" in column D.
$ws.Range("A81").Value = 20220413
$ws.Range("B81").Value = 2223.6496216999999
$ws.Range("C81").Value = 2224.4699999999998
$ws.Range("D81").Formula = "=100*(B81-C81)/C81"
$ws.Range("E81").Value = 180
$ws.Range("F81").Value = "CRM OPENED 20220401 MG"

# Update the selection to reflect where the user left off entering data.
$ws.Range("E82").Select()
